$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp in A1
$ws.Range("A1").Value = "Datos actualizados a 14 de Octubre de 2020 a las 17:37"

# Refresh country stats: the data feed re-sorted rows by "Casos totales"
# (desc), so several countries swapped rank/rows, and all numeric columns
# (B:H) received updated counts. Write each changed row in place.

$ws.Cells.Item(4, 1).Value = "Estados Unidos"
$ws.Cells.Item(4, 2).Value = 8103436
$ws.Cells.Item(4, 3).Value = 13183
$ws.Cells.Item(4, 4).Value = 5238565
$ws.Cells.Item(4, 5).Value = 2643769
$ws.Cells.Item(4, 6).Value = 0
$ws.Cells.Item(4, 7).Value = 229
$ws.Cells.Item(4, 8).Value = 221102

$ws.Cells.Item(5, 1).Value = "India"
$ws.Cells.Item(5, 2).Value = 7275588
$ws.Cells.Item(5, 3).Value = 38506
$ws.Cells.Item(5, 4).Value = 6343270
$ws.Cells.Item(5, 5).Value = 821352
$ws.Cells.Item(5, 6).Value = 0
$ws.Cells.Item(5, 7).Value = 349
$ws.Cells.Item(5, 8).Value = 110966

$ws.Cells.Item(15, 1).Value = "Reino Unido"
$ws.Cells.Item(15, 2).Value = 654644
$ws.Cells.Item(15, 3).Value = 19724
$ws.Cells.Item(15, 4).Value = 0
$ws.Cells.Item(15, 5).Value = 0
$ws.Cells.Item(15, 6).Value = 0
$ws.Cells.Item(15, 7).Value = 137
$ws.Cells.Item(15, 8).Value = 43155

$ws.Cells.Item(17, 1).Value = "Chile"
$ws.Cells.Item(17, 2).Value = 485372
$ws.Cells.Item(17, 3).Value = 1092
$ws.Cells.Item(17, 4).Value = 458073
$ws.Cells.Item(17, 5).Value = 13884
$ws.Cells.Item(17, 6).Value = 0
$ws.Cells.Item(17, 7).Value = 19
$ws.Cells.Item(17, 8).Value = 13415

$ws.Cells.Item(20, 1).Value = "Italia"
$ws.Cells.Item(20, 2).Value = 372799
$ws.Cells.Item(20, 3).Value = 7332
$ws.Cells.Item(20, 4).Value = 244065
$ws.Cells.Item(20, 5).Value = 92445
$ws.Cells.Item(20, 6).Value = 0
$ws.Cells.Item(20, 7).Value = 43
$ws.Cells.Item(20, 8).Value = 36289

$ws.Cells.Item(24, 1).Value = "Alemania"
$ws.Cells.Item(24, 2).Value = 339281
$ws.Cells.Item(24, 3).Value = 3602
$ws.Cells.Item(24, 4).Value = 281900
$ws.Cells.Item(24, 5).Value = 47623
$ws.Cells.Item(24, 6).Value = 0
$ws.Cells.Item(24, 7).Value = 18
$ws.Cells.Item(24, 8).Value = 9758

$ws.Cells.Item(25, 1).Value = "Turquia"
$ws.Cells.Item(25, 2).Value = 338779
$ws.Cells.Item(25, 3).Value = 0
$ws.Cells.Item(25, 4).Value = 296972
$ws.Cells.Item(25, 5).Value = 32850
$ws.Cells.Item(25, 6).Value = 0
$ws.Cells.Item(25, 7).Value = 0
$ws.Cells.Item(25, 8).Value = 8957

$ws.Cells.Item(30, 1).Value = "Canada"
$ws.Cells.Item(30, 2).Value = 187602
$ws.Cells.Item(30, 3).Value = 721
$ws.Cells.Item(30, 4).Value = 158269
$ws.Cells.Item(30, 5).Value = 19679
$ws.Cells.Item(30, 6).Value = 0
$ws.Cells.Item(30, 7).Value = 0
$ws.Cells.Item(30, 8).Value = 9654

$ws.Cells.Item(48, 1).Value = "Guatemala"
$ws.Cells.Item(48, 2).Value = 99094
$ws.Cells.Item(48, 3).Value = 714
$ws.Cells.Item(48, 4).Value = 88416
$ws.Cells.Item(48, 5).Value = 7248
$ws.Cells.Item(48, 6).Value = 0
$ws.Cells.Item(48, 7).Value = 20
$ws.Cells.Item(48, 8).Value = 3430

$ws.Cells.Item(49, 1).Value = "Portugal"
$ws.Cells.Item(49, 2).Value = 91193
$ws.Cells.Item(49, 3).Value = 2072
$ws.Cells.Item(49, 4).Value = 54493
$ws.Cells.Item(49, 5).Value = 34583
$ws.Cells.Item(49, 6).Value = 0
$ws.Cells.Item(49, 7).Value = 7
$ws.Cells.Item(49, 8).Value = 2117

$ws.Cells.Item(50, 1).Value = "Costa Rica"
$ws.Cells.Item(50, 2).Value = 90238
$ws.Cells.Item(50, 3).Value = 0
$ws.Cells.Item(50, 4).Value = 53670
$ws.Cells.Item(50, 5).Value = 35444
$ws.Cells.Item(50, 6).Value = 0
$ws.Cells.Item(50, 7).Value = 0
$ws.Cells.Item(50, 8).Value = 1124

$ws.Cells.Item(51, 1).Value = "Japon"
$ws.Cells.Item(51, 2).Value = 90140
$ws.Cells.Item(51, 3).Value = 467
$ws.Cells.Item(51, 4).Value = 83158
$ws.Cells.Item(51, 5).Value = 5344
$ws.Cells.Item(51, 6).Value = 0
$ws.Cells.Item(51, 7).Value = 4
$ws.Cells.Item(51, 8).Value = 1638

$ws.Cells.Item(59, 1).Value = "Moldavia"
$ws.Cells.Item(59, 2).Value = 64424
$ws.Cells.Item(59, 3).Value = 1149
$ws.Cells.Item(59, 4).Value = 45545
$ws.Cells.Item(59, 5).Value = 17365
$ws.Cells.Item(59, 6).Value = 0
$ws.Cells.Item(59, 7).Value = 19
$ws.Cells.Item(59, 8).Value = 1514

$ws.Cells.Item(64, 1).Value = "Singapur"
$ws.Cells.Item(64, 2).Value = 57889
$ws.Cells.Item(64, 3).Value = 5
$ws.Cells.Item(64, 4).Value = 57752
$ws.Cells.Item(64, 5).Value = 109
$ws.Cells.Item(64, 6).Value = 0
$ws.Cells.Item(64, 7).Value = 0
$ws.Cells.Item(64, 8).Value = 28

$ws.Cells.Item(81, 1).Value = "Birmania"
$ws.Cells.Item(81, 2).Value = 31325
$ws.Cells.Item(81, 3).Value = 888
$ws.Cells.Item(81, 4).Value = 13866
$ws.Cells.Item(81, 5).Value = 16727
$ws.Cells.Item(81, 6).Value = 0
$ws.Cells.Item(81, 7).Value = 39
$ws.Cells.Item(81, 8).Value = 732

$ws.Cells.Item(82, 1).Value = "El Salvador"
$ws.Cells.Item(82, 2).Value = 30766
$ws.Cells.Item(82, 3).Value = 286
$ws.Cells.Item(82, 4).Value = 26087
$ws.Cells.Item(82, 5).Value = 3775
$ws.Cells.Item(82, 6).Value = 0
$ws.Cells.Item(82, 7).Value = 5
$ws.Cells.Item(82, 8).Value = 904

$ws.Cells.Item(83, 1).Value = "Jordania"
$ws.Cells.Item(83, 2).Value = 30550
$ws.Cells.Item(83, 3).Value = 2423
$ws.Cells.Item(83, 4).Value = 6466
$ws.Cells.Item(83, 5).Value = 23827
$ws.Cells.Item(83, 6).Value = 0
$ws.Cells.Item(83, 7).Value = 32
$ws.Cells.Item(83, 8).Value = 257

$ws.Cells.Item(95, 1).Value = "Albania"
$ws.Cells.Item(95, 2).Value = 15955
$ws.Cells.Item(95, 3).Value = 203
$ws.Cells.Item(95, 4).Value = 9762
$ws.Cells.Item(95, 5).Value = 5759
$ws.Cells.Item(95, 6).Value = 0
$ws.Cells.Item(95, 7).Value = 5
$ws.Cells.Item(95, 8).Value = 434

$ws.Cells.Item(96, 1).Value = "Noruega"
$ws.Cells.Item(96, 2).Value = 15888
$ws.Cells.Item(96, 3).Value = 97
$ws.Cells.Item(96, 4).Value = 11863
$ws.Cells.Item(96, 5).Value = 3748
$ws.Cells.Item(96, 6).Value = 0
$ws.Cells.Item(96, 7).Value = 0
$ws.Cells.Item(96, 8).Value = 277

$ws.Cells.Item(99, 1).Value = "Montenegro"
$ws.Cells.Item(99, 2).Value = 14461
$ws.Cells.Item(99, 3).Value = 193
$ws.Cells.Item(99, 4).Value = 10201
$ws.Cells.Item(99, 5).Value = 4043
$ws.Cells.Item(99, 6).Value = 0
$ws.Cells.Item(99, 7).Value = 4
$ws.Cells.Item(99, 8).Value = 217

$ws.Cells.Item(111, 1).Value = "Luxemburgo"
$ws.Cells.Item(111, 2).Value = 10030
$ws.Cells.Item(111, 3).Value = 190
$ws.Cells.Item(111, 4).Value = 8306
$ws.Cells.Item(111, 5).Value = 1591
$ws.Cells.Item(111, 6).Value = 0
$ws.Cells.Item(111, 7).Value = 0
$ws.Cells.Item(111, 8).Value = 133

$ws.Cells.Item(112, 1).Value = "Eslovenia"
$ws.Cells.Item(112, 2).Value = 9938
$ws.Cells.Item(112, 3).Value = 707
$ws.Cells.Item(112, 4).Value = 5515
$ws.Cells.Item(112, 5).Value = 4248
$ws.Cells.Item(112, 6).Value = 0
$ws.Cells.Item(112, 7).Value = 2
$ws.Cells.Item(112, 8).Value = 175

$ws.Cells.Item(116, 1).Value = "Jamaica"
$ws.Cells.Item(116, 2).Value = 7989
$ws.Cells.Item(116, 3).Value = 79
$ws.Cells.Item(116, 4).Value = 3431
$ws.Cells.Item(116, 5).Value = 4407
$ws.Cells.Item(116, 6).Value = 0
$ws.Cells.Item(116, 7).Value = 5
$ws.Cells.Item(116, 8).Value = 151

$ws.Cells.Item(122, 1).Value = "Cuba"
$ws.Cells.Item(122, 2).Value = 6035
$ws.Cells.Item(122, 3).Value = 18
$ws.Cells.Item(122, 4).Value = 5653
$ws.Cells.Item(122, 5).Value = 259
$ws.Cells.Item(122, 6).Value = 0
$ws.Cells.Item(122, 7).Value = 0
$ws.Cells.Item(122, 8).Value = 123

$ws.Cells.Item(128, 1).Value = "Sri Lanka"
$ws.Cells.Item(128, 2).Value = 5168
$ws.Cells.Item(128, 3).Value = 130
$ws.Cells.Item(128, 4).Value = 3357
$ws.Cells.Item(128, 5).Value = 1798
$ws.Cells.Item(128, 6).Value = 0
$ws.Cells.Item(128, 7).Value = 0
$ws.Cells.Item(128, 8).Value = 13

$ws.Cells.Item(129, 1).Value = "Bahamas"
$ws.Cells.Item(129, 2).Value = 5163
$ws.Cells.Item(129, 3).Value = 0
$ws.Cells.Item(129, 4).Value = 2978
$ws.Cells.Item(129, 5).Value = 2077
$ws.Cells.Item(129, 6).Value = 0
$ws.Cells.Item(129, 7).Value = 0
$ws.Cells.Item(129, 8).Value = 108

$ws.Cells.Item(130, 1).Value = "Trinidad yTobago"
$ws.Cells.Item(130, 2).Value = 5127
$ws.Cells.Item(130, 3).Value = 0
$ws.Cells.Item(130, 4).Value = 3367
$ws.Cells.Item(130, 5).Value = 1667
$ws.Cells.Item(130, 6).Value = 0
$ws.Cells.Item(130, 7).Value = 0
$ws.Cells.Item(130, 8).Value = 93

$ws.Cells.Item(131, 1).Value = "Congo"
$ws.Cells.Item(131, 2).Value = 5118
$ws.Cells.Item(131, 3).Value = 0
$ws.Cells.Item(131, 4).Value = 3887
$ws.Cells.Item(131, 5).Value = 1141
$ws.Cells.Item(131, 6).Value = 0
$ws.Cells.Item(131, 7).Value = 0
$ws.Cells.Item(131, 8).Value = 90

$ws.Cells.Item(132, 1).Value = "Surinam"
$ws.Cells.Item(132, 2).Value = 5072
$ws.Cells.Item(132, 3).Value = 0
$ws.Cells.Item(132, 4).Value = 4870
$ws.Cells.Item(132, 5).Value = 95
$ws.Cells.Item(132, 6).Value = 0
$ws.Cells.Item(132, 7).Value = 0
$ws.Cells.Item(132, 8).Value = 107

$ws.Cells.Item(133, 1).Value = "Guinea Ecuatorial"
$ws.Cells.Item(133, 2).Value = 5066
$ws.Cells.Item(133, 3).Value = 0
$ws.Cells.Item(133, 4).Value = 4954
$ws.Cells.Item(133, 5).Value = 29
$ws.Cells.Item(133, 6).Value = 0
$ws.Cells.Item(133, 7).Value = 0
$ws.Cells.Item(133, 8).Value = 83

$ws.Cells.Item(137, 1).Value = "Reunion"
$ws.Cells.Item(137, 2).Value = 4678
$ws.Cells.Item(137, 3).Value = 54
$ws.Cells.Item(137, 4).Value = 3994
$ws.Cells.Item(137, 5).Value = 667
$ws.Cells.Item(137, 6).Value = 0
$ws.Cells.Item(137, 7).Value = 0
$ws.Cells.Item(137, 8).Value = 17

$ws.Cells.Item(174, 1).Value = "Curazao"
$ws.Cells.Item(174, 2).Value = 645
$ws.Cells.Item(174, 3).Value = 26
$ws.Cells.Item(174, 4).Value = 367
$ws.Cells.Item(174, 5).Value = 277
$ws.Cells.Item(174, 6).Value = 0
$ws.Cells.Item(174, 7).Value = 0
$ws.Cells.Item(174, 8).Value = 1

$ws.Cells.Item(175, 1).Value = "Papua Nueva Guinea"
$ws.Cells.Item(175, 2).Value = 578
$ws.Cells.Item(175, 3).Value = 13
$ws.Cells.Item(175, 4).Value = 537
$ws.Cells.Item(175, 5).Value = 34
$ws.Cells.Item(175, 6).Value = 0
$ws.Cells.Item(175, 7).Value = 0
$ws.Cells.Item(175, 8).Value = 7

$ws.Cells.Item(197, 1).Value = "Antigua y Barbuda"
$ws.Cells.Item(197, 2).Value = 112
$ws.Cells.Item(197, 3).Value = 1
$ws.Cells.Item(197, 4).Value = 100
$ws.Cells.Item(197, 5).Value = 9
$ws.Cells.Item(197, 6).Value = 0
$ws.Cells.Item(197, 7).Value = 0
$ws.Cells.Item(197, 8).Value = 3

